$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 20240112
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 6

$ws.Range("B8").Value = "24 trials-->"

$ws.Range("A9").Value = 20240119
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 4

$ws.Range("A10").Select()
